$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1714

# Row 3
$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 1143

# Row 4
$ws.Range("D4").Value = 44299
$ws.Range("M4").Value = 80
$ws.Range("R4").Value = "Provincia de Santiago"

# Row 5
$ws.Range("D5").Value = 44299
$ws.Range("M5").Value = 75
$ws.Range("R5").Value = "Provincia de Santiago"

# Row 6
$ws.Range("D6").Value = 44292
$ws.Range("M6").Value = 25
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("S6").Value = 2286

# Row 7
$ws.Range("D7").Value = 44292
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 2143

# Row 8
$ws.Range("D8").Value = 44302
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 2143

# Row 9
$ws.Range("D9").Value = 44302
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 10
$ws.Range("D10").Value = 44322
$ws.Range("M10").Value = 45

# Row 11
$ws.Range("D11").Value = 44322
$ws.Range("M11").Value = 80

# Row 12
$ws.Range("D12").Value = 44300
$ws.Range("M12").Value = 100

# Row 13
$ws.Range("D13").Value = 44300
$ws.Range("M13").Value = 80

# Row 14
$ws.Range("D14").Value = 44301
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 14000
$ws.Range("S14").Value = 2000

# Row 15
$ws.Range("D15").Value = 44301
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 12000
$ws.Range("S15").Value = 1714
